$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute("testdate", $true, $false, $false, $false, $false, $true, 1, $false, "test_date", 2)
